# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Goblin Profits leve-crafting tables
# across all eight sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 852.25
$ws.Range("I20").Value = 869.6667
$ws.Range("J20").Value = 800
$ws.Range("K20").Value = 869.6667
$ws.Range("L20").Value = 800
$ws.Range("M20").Value = -639.6667
$ws.Range("N20").Value = -1260
$ws.Range("H32").Value = 10000
$ws.Range("I32").Value = 10000
$ws.Range("K32").Value = 10000
$ws.Range("M32").Value = -9674
$ws.Range("H35").Value = 852.25
$ws.Range("I35").Value = 869.6667
$ws.Range("J35").Value = 800
$ws.Range("K35").Value = 869.6667
$ws.Range("L35").Value = 800
$ws.Range("M35").Value = -490.6667
$ws.Range("N35").Value = -1558
$ws.Range("H38").Value = 7121.7036
$ws.Range("I38").Value = 7393.25
$ws.Range("J38").Value = 6726.727
$ws.Range("K38").Value = 22179.75
$ws.Range("L38").Value = 20180.181
$ws.Range("M38").Value = -21807.75
$ws.Range("N38").Value = -20924.181
$ws.Range("H69").Value = 111115450
$ws.Range("I69").Value = 3999.6667
$ws.Range("J69").Value = 166671170
$ws.Range("K69").Value = 11999.0001
$ws.Range("L69").Value = 500013510
$ws.Range("M69").Value = -11125.0001
$ws.Range("N69").Value = -500015258
$ws.Range("H72").Value = 111115450
$ws.Range("I72").Value = 3999.6667
$ws.Range("J72").Value = 166671170
$ws.Range("K72").Value = 35997.0003
$ws.Range("L72").Value = 1500040530
$ws.Range("M72").Value = -31629.0003
$ws.Range("N72").Value = -1500049266
$ws.Range("H96").Value = 667544.25
$ws.Range("I96").Value = 1111599.5
$ws.Range("K96").Value = 3334798.5
$ws.Range("M96").Value = -3333425.5
$ws.Range("H100").Value = 6374.294
$ws.Range("I100").Value = 4107
$ws.Range("K100").Value = 4107
$ws.Range("M100").Value = -3566
$ws.Range("H112").Value = 1793.2858
$ws.Range("J112").Value = 2345.7144
$ws.Range("L112").Value = 7037.1432
$ws.Range("N112").Value = -9253.143199999999
$ws.Range("H132").Value = 1380.6171
$ws.Range("I132").Value = 1199.079
$ws.Range("K132").Value = 3597.237
$ws.Range("M132").Value = -1067.237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2922.9
$ws.Range("I32").Value = 2790.5532
$ws.Range("K32").Value = 2790.5532
$ws.Range("M32").Value = -2503.5532
$ws.Range("H63").Value = 3893
$ws.Range("J63").Value = 2474.5
$ws.Range("L63").Value = 2474.5
$ws.Range("N63").Value = -3846.5
$ws.Range("H66").Value = 3893
$ws.Range("J66").Value = 2474.5
$ws.Range("L66").Value = 12372.5
$ws.Range("N66").Value = -19236.5
$ws.Range("H97").Value = 2236.6155
$ws.Range("I97").Value = 2029.6364
$ws.Range("J97").Value = 3375
$ws.Range("K97").Value = 2029.6364
$ws.Range("L97").Value = 3375
$ws.Range("M97").Value = -1533.6364
$ws.Range("N97").Value = -4367

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 6897.8
$ws.Range("I33").Value = 7874.75
$ws.Range("J33").Value = 2990
$ws.Range("K33").Value = 7874.75
$ws.Range("L33").Value = 2990
$ws.Range("M33").Value = -7538.75
$ws.Range("N33").Value = -3662
$ws.Range("H86").Value = 577576.8
$ws.Range("I86").Value = 2213.25
$ws.Range("K86").Value = 2213.25
$ws.Range("M86").Value = -1090.25
$ws.Range("H89").Value = 577576.8
$ws.Range("I89").Value = 2213.25
$ws.Range("K89").Value = 11066.25
$ws.Range("M89").Value = -5450.25
$ws.Range("H94").Value = 3356.25
$ws.Range("I94").Value = 3222
$ws.Range("K94").Value = 3222
$ws.Range("M94").Value = -2771
$ws.Range("H99").Value = 3649.8
$ws.Range("J99").Value = 4333.1665
$ws.Range("L99").Value = 4333.1665
$ws.Range("N99").Value = -7329.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2980.875
$ws.Range("I31").Value = 2441.3333
$ws.Range("J31").Value = 4599.5
$ws.Range("K31").Value = 2441.3333
$ws.Range("L31").Value = 4599.5
$ws.Range("M31").Value = -2146.3333
$ws.Range("N31").Value = -5189.5
$ws.Range("H34").Value = 2980.875
$ws.Range("I34").Value = 2441.3333
$ws.Range("J34").Value = 4599.5
$ws.Range("K34").Value = 2441.3333
$ws.Range("L34").Value = 4599.5
$ws.Range("M34").Value = -2239.3333
$ws.Range("N34").Value = -5003.5
$ws.Range("H75").Value = 19399.6
$ws.Range("J75").Value = 19399.6
$ws.Range("L75").Value = 19399.6
$ws.Range("N75").Value = -21395.6
$ws.Range("H78").Value = 19399.6
$ws.Range("J78").Value = 19399.6
$ws.Range("L78").Value = 58198.8
$ws.Range("N78").Value = -68182.79999999999
$ws.Range("H99").Value = 4025
$ws.Range("I99").Value = 4025
$ws.Range("K99").Value = 4025
$ws.Range("M99").Value = -2527
$ws.Range("H126").Value = 4025
$ws.Range("I126").Value = 4025
$ws.Range("K126").Value = 12075
$ws.Range("M126").Value = -9605
$ws.Range("H132").Value = 4465.9287
$ws.Range("I132").Value = 4655.6924
$ws.Range("K132").Value = 13967.0772
$ws.Range("M132").Value = -11437.0772
$ws.Range("H134").Value = 4681.364
$ws.Range("J134").Value = 7000
$ws.Range("L134").Value = 21000
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2114.6667
$ws.Range("I92").Value = 1671.25
$ws.Range("J92").Value = 3001.5
$ws.Range("K92").Value = 5013.75
$ws.Range("L92").Value = 9004.5
$ws.Range("M92").Value = -3765.75
$ws.Range("N92").Value = -11500.5
$ws.Range("H93").Value = 13758.167
$ws.Range("I93").Value = 859.2
$ws.Range("K93").Value = 2577.6
$ws.Range("M93").Value = -705.6000000000004
$ws.Range("H96").Value = 3000
$ws.Range("J96").Value = 3000
$ws.Range("L96").Value = 9000
$ws.Range("N96").Value = -13118
$ws.Range("H98").Value = 913.8570999999999
$ws.Range("I98").Value = 799.3333
$ws.Range("K98").Value = 2397.9999
$ws.Range("M98").Value = -899.9998999999998
$ws.Range("H99").Value = 32010
$ws.Range("I99").Value = 16683.334
$ws.Range("K99").Value = 50050.00199999999
$ws.Range("M99").Value = -47804.00199999999
$ws.Range("H100").Value = 13303.167
$ws.Range("J100").Value = 34850
$ws.Range("L100").Value = 104550
$ws.Range("N100").Value = -106172
$ws.Range("H132").Value = 2012.5454
$ws.Range("I132").Value = 1706.8462
$ws.Range("J132").Value = 2454.111
$ws.Range("K132").Value = 15361.6158
$ws.Range("L132").Value = 22086.999
$ws.Range("M132").Value = -12831.6158
$ws.Range("N132").Value = -27146.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5129.696
$ws.Range("I80").Value = 4165.6
$ws.Range("J80").Value = 6937.375
$ws.Range("K80").Value = 4165.6
$ws.Range("L80").Value = 6937.375
$ws.Range("M80").Value = -3167.6
$ws.Range("N80").Value = -8933.375
$ws.Range("H83").Value = 5129.696
$ws.Range("I83").Value = 4165.6
$ws.Range("J83").Value = 6937.375
$ws.Range("K83").Value = 20828
$ws.Range("L83").Value = 34686.875
$ws.Range("M83").Value = -15836
$ws.Range("N83").Value = -44670.875
$ws.Range("H97").Value = 1081.5454
$ws.Range("I97").Value = 916.6667
$ws.Range("K97").Value = 916.6667
$ws.Range("M97").Value = -420.6667
$ws.Range("H104").Value = 33690
$ws.Range("J104").Value = 33690
$ws.Range("L104").Value = 33690
$ws.Range("N104").Value = -40678
$ws.Range("H122").Value = 10030.037
$ws.Range("I122").Value = 11719
$ws.Range("K122").Value = 35157
$ws.Range("M122").Value = -32707
$ws.Range("H126").Value = 2619.4736
$ws.Range("I126").Value = 2698.7144
$ws.Range("J126").Value = 2397.6
$ws.Range("K126").Value = 8096.1432
$ws.Range("L126").Value = 7192.799999999999
$ws.Range("M126").Value = -5626.1432
$ws.Range("N126").Value = -12132.8
$ws.Range("H132").Value = 2786.4
$ws.Range("I132").Value = 2716
$ws.Range("K132").Value = 8148
$ws.Range("M132").Value = -5618

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2866.3572
$ws.Range("I68").Value = 2248.0908
$ws.Range("J68").Value = 5133.3335
$ws.Range("K68").Value = 2248.0908
$ws.Range("L68").Value = 5133.3335
$ws.Range("M68").Value = -1499.0908
$ws.Range("N68").Value = -6631.3335
$ws.Range("H71").Value = 2866.3572
$ws.Range("I71").Value = 2248.0908
$ws.Range("J71").Value = 5133.3335
$ws.Range("K71").Value = 11240.454
$ws.Range("L71").Value = 25666.6675
$ws.Range("M71").Value = -7496.454
$ws.Range("N71").Value = -33154.6675
$ws.Range("H123").Value = 113883.664
$ws.Range("J123").Value = 150825.5
$ws.Range("L123").Value = 150825.5
$ws.Range("N123").Value = -160625.5
$ws.Range("H136").Value = 8497.5
$ws.Range("J136").Value = 8497.5
$ws.Range("L136").Value = 25492.5
$ws.Range("N136").Value = -30592.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3102.1052
$ws.Range("I81").Value = 2997.2307
$ws.Range("J81").Value = 3329.3333
$ws.Range("K81").Value = 5994.4614
$ws.Range("L81").Value = 6658.6666
$ws.Range("M81").Value = -4933.4614
$ws.Range("N81").Value = -8780.6666
$ws.Range("H84").Value = 3102.1052
$ws.Range("I84").Value = 2997.2307
$ws.Range("J84").Value = 3329.3333
$ws.Range("K84").Value = 29972.307
$ws.Range("L84").Value = 33293.333
$ws.Range("M84").Value = -24668.307
$ws.Range("N84").Value = -43901.333
$ws.Range("H96").Value = 3817.6
$ws.Range("J96").Value = 4580.5
$ws.Range("L96").Value = 4580.5
$ws.Range("N96").Value = -7326.5
$ws.Range("H122").Value = 4749.385
$ws.Range("I122").Value = 4749.385
$ws.Range("K122").Value = 14248.155
$ws.Range("M122").Value = -11798.155
$ws.Range("H132").Value = 1618.0526
$ws.Range("I132").Value = 1618.0526
$ws.Range("K132").Value = 4854.1578
$ws.Range("M132").Value = -2324.1578
